$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the category labels in column A (Insert -> Ins, Delete -> Del)
$ws.Range("A2").Value = "Ins 100"
$ws.Range("A3").Value = "Ins 1,000"
$ws.Range("A4").Value = "Ins 10,000"
$ws.Range("A5").Value = "Ins 100,000"
$ws.Range("A6").Value = "Del 100"
$ws.Range("A7").Value = "Del 1,000"
$ws.Range("A9").Value = "Del 100,000"
$ws.Range("A8").Value = "Del 10,000"

# Update the active cell selection to A2
$ws.Range("A2").Select()
